$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Franco") deliverable answers
$ws.Range("E5").Value  = "si"
$ws.Range("E6").Value  = "no"
$ws.Range("E7").Value  = "si"
$ws.Range("E8").Value  = "si(lo tengo que sacar)"
$ws.Range("E9").Value  = "si"
$ws.Range("E10").Value = "si"
$ws.Range("E11").Value = "si"
$ws.Range("E12").Value = "si"
$ws.Range("E13").Value = "si"
$ws.Range("E14").Value = "si"
$ws.Range("E15").Value = "si"
$ws.Range("E16").Value = "si"
$ws.Range("E17").Value = "si"
$ws.Range("E18").Value = "no"
$ws.Range("E19").Value = "si"
$ws.Range("E20").Value = "si"
$ws.Range("E21").Value = "si"
$ws.Range("E22").Value = "si"
$ws.Range("E23").Value = "si"
$ws.Range("E25").Value = "si"
$ws.Range("E26").Value = "no"
$ws.Range("E27").Value = "si(tambien tengo gestion de riesgos)"
$ws.Range("E28").Value = "si"

# Extra comment note, shown in red font
$ws.Range("E30").Value = "tengo un par de entregables mas. Lo que me hace ruido es que algunos de mis entregables son las hojas del arbol, y algunas son los nodos. Como seria??"
$ws.Range("E30").Font.Color = 255

# Selection/view tweaks
$ws.Range("D12").Select() | Out-Null
$ws.Application.ActiveWindow.Zoom = 100
